# Book1.xlsx edit: update RUNMANAGER!C6, move both sheets' selections,
# and make DATA the active/selected sheet (as captured by the workbook's
# last-saved view state).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("DATA")

# RUNMANAGER!C6 "yes" -> "no"
$ws1.Range("C6").Value = "no"

# RUNMANAGER's own selection moves to A2 (captured while it's the active
# sheet, before we switch away).
$ws1.Activate()
$ws1.Range("A2").Select()

# DATA becomes the active/selected sheet with its selection on B4.
$ws2.Activate()
$ws2.Range("B4").Select()
